$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.783.79"
$ws.Range("E2").Value = "  +0.60%  "

$ws.Range("D3").Value = "1.639.18"
$ws.Range("E3").Value = "  +0.54%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.997"
$ws.Range("E4").Value = "  -0.50%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "212.60"
$ws.Range("E5").Value = "  +0.28%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.523"
$ws.Range("E6").Value = "  +0.23%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.998"
$ws.Range("E7").Value = "  -0.39%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "23.54"
$ws.Range("E8").Value = "  +1.89%  "

$ws.Range("E9").Value = "  -1.92%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0614"
$ws.Range("E10").Value = "  +0.35%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0881"
$ws.Range("E11").Value = "  +2.26%  "

$ws.Range("D12").Value = "1.872.37"
$ws.Range("E12").Value = "  +0.47%  "

$ws.Range("D13").Value = "1.640.07"
$ws.Range("E13").Value = "  +0.50%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.574"

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.09"
$ws.Range("E15").Value = "  +1.32%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.78"
$ws.Range("E16").Value = "  +1.13%  "

$ws.Range("D17").Value = "27.833.01"
$ws.Range("E17").Value = "  +0.88%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "233.33"
$ws.Range("E18").Value = "  +1.92%  "

$ws.Range("D19").Value = "0.0₃0723"
$ws.Range("E19").Value = "  +0.54%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.60"
$ws.Range("E20").Value = "  +0.78%  "

$ws.Range("E21").Value = "  -0.27%  "

$ws.Range("B22").Value = "Uniswap"
$ws.Range("C22").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.37"
$ws.Range("E22").Value = "  +0.54%  "

$ws.Range("B23").Value = "Avalanche"
$ws.Range("C23").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.37"
$ws.Range("E23").Value = "  -2.57%  "

$ws.Range("E24").Value = "  -1.69%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "152.32"
$ws.Range("E25").Value = "  +2.22%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.93"
$ws.Range("E26").Value = "  +0.85%  "

$ws.Range("E28").Value = "  +0.24%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.999"
$ws.Range("E29").Value = "  -0.40%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.18"
$ws.Range("E30").Value = "  -0.03%  "

$ws.Range("E31").Value = "  +0.66%  "

$ws.Range("E32").Value = "  +2.54%  "

$ws.Range("E33").Value = "  +1.40%  "

$ws.Range("D34").Value = "1.414.92"
$ws.Range("E34").Value = "  -3.81%  "

$ws.Range("E35").Value = "  +2.46%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.35"
$ws.Range("E36").Value = "  +1.11%  "

$ws.Range("E37").Value = "  +1.46%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.881"
$ws.Range("E38").Value = "  +0.50%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.560"
$ws.Range("E39").Value = "  +0.44%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.906"
$ws.Range("E40").Value = "  -2.02%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.02"
$ws.Range("E41").Value = "  +0.83%  "

$ws.Range("E42").Value = "  -0.32%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.87"
$ws.Range("E43").Value = "  +7.35%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "66.99"
$ws.Range("E44").Value = "  -1.20%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "5.52"
$ws.Range("E45").Value = "  +3.22%  "

$ws.Range("E46").Value = "  -0.15%  "

$ws.Range("D47").Value = "1.781.89"
$ws.Range("E47").Value = "  +0.54%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "87.63"
$ws.Range("E48").Value = "  +0.32%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0999"
$ws.Range("E49").Value = "  +0.79%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0505"
$ws.Range("E50").Value = "  +0.25%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.58"
$ws.Range("E51").Value = "  -1.17%  "
